# Fruta / hortaliza, semanal
# Insert a new weekly record as row 93 on the sheet, shifting the existing
# rows 93-108 down to 94-109 (the last old row becomes the new row 109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93; Excel shifts rows 93:108 down to 94:109
# and copies formatting (including the date style on column D) from the
# row above, matching the rest of the table.
$ws.Rows("93:93").Insert()

$ws.Range("A93").Value = 11
$ws.Range("B93").Value = "Vega Monumental Concepción"
$ws.Range("C93").Value = "Bíobío"
$ws.Range("D93").Value = 44637
$ws.Range("E93").Value = 8
$ws.Range("F93").Value = 100112032
$ws.Range("G93").Value = "Zapallo italiano"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 220
$ws.Range("K93").Value = 12000
$ws.Range("L93").Value = 14000
$ws.Range("M93").Value = 12909
$ws.Range("N93").Value = "$/caja 50 unidades"
$ws.Range("O93").Value = "Región de O'Higgins"
$ws.Range("P93").Value = 258
$ws.Range("Q93").Value = 50
$ws.Range("R93").Value = "Hortaliza"
